$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 486, shifting existing rows 486-553 down to 488-555
$ws.Rows("486:487").Insert()

# Row 486
$ws.Cells.Item(486,1).Value = 5
$ws.Cells.Item(486,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(486,3).Value = 'Maule'
$ws.Cells.Item(486,4).Value = 44776
$ws.Cells.Item(486,5).Value = 7
$ws.Cells.Item(486,6).Value = 100114001
$ws.Cells.Item(486,7).Value = 'Papa'
$ws.Cells.Item(486,8).Value = 'Rodeo'
$ws.Cells.Item(486,9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(486,10).Value = 1200
$ws.Cells.Item(486,11).Value = 8000
$ws.Cells.Item(486,12).Value = 8000
$ws.Cells.Item(486,13).Value = 8000
$ws.Cells.Item(486,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(486,15).Value = 'Región de La Araucanía'
$ws.Cells.Item(486,16).Value = 320
$ws.Cells.Item(486,17).Value = 25
$ws.Cells.Item(486,18).Value = 'Hortaliza'
$ws.Cells.Item(486,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 487
$ws.Cells.Item(487,1).Value = 5
$ws.Cells.Item(487,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(487,3).Value = 'Maule'
$ws.Cells.Item(487,4).Value = 44776
$ws.Cells.Item(487,5).Value = 7
$ws.Cells.Item(487,6).Value = 100114001
$ws.Cells.Item(487,7).Value = 'Papa'
$ws.Cells.Item(487,8).Value = 'Rosara'
$ws.Cells.Item(487,9).Value = '1a (cosecha)'
$ws.Cells.Item(487,10).Value = 1500
$ws.Cells.Item(487,11).Value = 5800
$ws.Cells.Item(487,12).Value = 5800
$ws.Cells.Item(487,13).Value = 5800
$ws.Cells.Item(487,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(487,15).Value = 'Región del Maule'
$ws.Cells.Item(487,16).Value = 232
$ws.Cells.Item(487,17).Value = 25
$ws.Cells.Item(487,18).Value = 'Hortaliza'
$ws.Cells.Item(487,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

